# Insert a new data row at row 58 (pushing existing rows 58..91 down to 59..92)
# and populate it with the new "Locoto" price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(58).Insert()

$ws.Range("A58").Value = 10
$ws.Range("B58").Value = "Vega Modelo de Temuco"
$ws.Range("C58").Value = "La Araucanía"
$ws.Range("D58").Value = 45126
$ws.Range("E58").Value = 9
$ws.Range("F58").Value = 100112042
$ws.Range("G58").Value = "Locoto"
$ws.Range("H58").Value = "Sin especificar"
$ws.Range("I58").Value = "Primera"
$ws.Range("J58").Value = 60
$ws.Range("K58").Value = 3800
$ws.Range("L58").Value = 3800
$ws.Range("M58").Value = 3800
$ws.Range("N58").Value = "$/kilo"
$ws.Range("O58").Value = "Región de Arica y Parinacota"
$ws.Range("P58").Value = 3800
$ws.Range("Q58").Value = 1
$ws.Range("R58").Value = "Hortaliza"

# Match the D-column (date) number format/style used by the other rows.
$ws.Range("D58").NumberFormat = $ws.Range("D59").NumberFormat
